$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capitalize hex letter-digits (a-f -> A-F) in the doip (G) and uds (H) columns
# for rows 2-33, e.g. "0xfd" -> "0xFD". The "0x" prefix itself stays lowercase.
$ws.Range("G2").Value = "0x02:0xFD:0x00:0x05:0x00:0x00:0x00:0x07:0x0E:0x00:0x00:0x00:0x00:0x00:0x00"
$ws.Range("G3").Value = "0x02:0xFD:0x00:0x06:0x00:0x00:0x00:0x0D:0x0E:0x00:0xE0:0x00:0x10:0x00:0x00:0x00:0x00:0x00:0x00:0x00:0x00"
$ws.Range("G4").Value = "0x02:0xFD:0x00:0x01:0x00:0x00:0x00:0x00"
$ws.Range("G5").Value = "0x02:0xFD:0x00:0x04:0x00:0x00:0x00:0x21:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0xE0:0x00:0xE1:0xE2:0xE3:0xE4:0xE5:0xE6:0xA1:0xA2:0xA3:0xA4:0xA5:0xA6:0x00:0x00"
$ws.Range("G6").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0x10:0x32:0x3E:0x00"
$ws.Range("H6").Value = "0x3E:0x00"
$ws.Range("G7").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0x10:0x32:0x0E:0x00:0x00"
$ws.Range("G8").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x10:0x32:0x0E:0x00:0x7E:0x00"
$ws.Range("H8").Value = "0x7E:0x00"
$ws.Range("G9").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0x10:0x32:0x3E:0x80"
$ws.Range("H9").Value = "0x3E:0x80"
$ws.Range("G10").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0x10:0x32:0x0E:0x00:0x00"
$ws.Range("G11").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x10:0x32:0x0E:0x00:0x7E:0x80"
$ws.Range("H11").Value = "0x7E:0x80"
$ws.Range("G12").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x07:0x0E:0x00:0x10:0x32:0x22:0xF1:0x90"
$ws.Range("H12").Value = "0x22:0xF1:0x90"
$ws.Range("G13").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0x10:0x32:0x0E:0x00:0x00"
$ws.Range("G14").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x18:0x10:0x32:0x0E:0x00:0x62:0xF1:0x90:0x4D:0x41:0x54:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34"
$ws.Range("H14").Value = "0x62:0xF1:0x90:0x4D:0x41:0x54:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34"
$ws.Range("G15").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0x10:0x32:0x10:0x01"
$ws.Range("G16").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0x10:0x32:0x0E:0x00:0x00"
$ws.Range("G17").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x0A:0x10:0x32:0x0E:0x00:0x50:0x01:0x00:0x32:0x01:0xF4"
$ws.Range("H17").Value = "0x50:0x01:0x00:0x32:0x01:0xF4"
$ws.Range("G18").Value = "0x02:0xFD:0x00:0x05:0x00:0x00:0x00:0x07:0x0E:0x00:0x00:0x00:0x00:0x00:0x00"
$ws.Range("G19").Value = "0x02:0xFD:0x00:0x06:0x00:0x00:0x00:0x0D:0x0E:0x00:0xE0:0x00:0x10:0x00:0x00:0x00:0x00:0x00:0x00:0x00:0x00"
$ws.Range("G20").Value = "0x02:0xFD:0x00:0x01:0x00:0x00:0x00:0x00"
$ws.Range("G21").Value = "0x02:0xFD:0x00:0x04:0x00:0x00:0x00:0x21:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0xE0:0x00:0xE1:0xE2:0xE3:0xE4:0xE5:0xE6:0xA1:0xA2:0xA3:0xA4:0xA5:0xA6:0x00:0x00"
$ws.Range("G22").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0xE0:0x01:0x3E:0x00"
$ws.Range("H22").Value = "0x3E:0x00"
$ws.Range("G23").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0xE0:0x01:0x0E:0x00:0x00"
$ws.Range("G24").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0xE0:0x01:0x0E:0x00:0x7E:0x00"
$ws.Range("H24").Value = "0x7E:0x00"
$ws.Range("G25").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0xE0:0x01:0x3E:0x80"
$ws.Range("H25").Value = "0x3E:0x80"
$ws.Range("G26").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0xE0:0x01:0x0E:0x00:0x00"
$ws.Range("G27").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0xE0:0x01:0x0E:0x00:0x7E:0x80"
$ws.Range("H27").Value = "0x7E:0x80"
$ws.Range("G28").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x07:0x0E:0x00:0xE0:0x01:0x22:0xF1:0x90"
$ws.Range("H28").Value = "0x22:0xF1:0x90"
$ws.Range("G29").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0xE0:0x01:0x0E:0x00:0x00"
$ws.Range("G30").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x18:0xE0:0x01:0x0E:0x00:0x62:0xF1:0x90:0x4D:0x41:0x54:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34"
$ws.Range("H30").Value = "0x62:0xF1:0x90:0x4D:0x41:0x54:0x31:0x32:0x33:0x34:0x35:0x36:0x37:0x38:0x39:0x30:0x31:0x32:0x33:0x34"
$ws.Range("G31").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x06:0x0E:0x00:0xE0:0x01:0x10:0x01"
$ws.Range("G32").Value = "0x02:0xFD:0x80:0x02:0x00:0x00:0x00:0x05:0xE0:0x01:0x0E:0x00:0x00"
$ws.Range("G33").Value = "0x02:0xFD:0x80:0x01:0x00:0x00:0x00:0x07:0xE0:0x01:0x0E:0x00:0x7F:0x10:0x13"
$ws.Range("H33").Value = "0x7F:0x10:0x13"
